$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 219: LeetCode 717, 1-bit and 2-bit Characters ---
$ws.Range("A219").Value2 = 717
$ws.Range("B219").Value2 = "1-bit and 2-bit Characters"
$ws.Range("C219").Value2 = "#array #math "
$ws.Range("D219").Value2 = "easy"
$ws.Range("E219").Value2 = 0
$ws.Range("F219").Value2 = 1
$ws.Range("G219").Value2 = 25
$ws.Range("H219").Value2 = 45980
$ws.Range("I219").Value2 = 45980
$ws.Rows.Item(219).RowHeight = 34

# --- Row 220: LeetCode 2154, Keep Multiplying Found Values by Two ---
$ws.Range("A220").Value2 = 2154
$ws.Range("B220").Value2 = "Keep Multiplying Found Values by Two"
$ws.Range("C220").Value2 = "#array"
$ws.Range("D220").Value2 = "easy"
$ws.Range("E220").Value2 = 0
$ws.Range("F220").Value2 = 1
$ws.Range("G220").Value2 = 5
$ws.Range("H220").Value2 = 45980
$ws.Range("I220").Value2 = 45980
$ws.Rows.Item(220).RowHeight = 34

# --- Row 221: LeetCode 757, Set Intersection Size At Least Two ---
$ws.Range("A221").Value2 = 757
$ws.Range("B221").Value2 = "Set Intersection Size At Least Two"
$ws.Range("C221").Value2 = "#array #intervals"
$ws.Range("D221").Value2 = "hard"
$ws.Range("E221").Value2 = 0
$ws.Range("F221").Value2 = 1
$ws.Range("G221").Value2 = 55
$ws.Range("H221").Value2 = 45981
$ws.Range("I221").Value2 = 45981
$ws.Rows.Item(221).RowHeight = 34

# Copy the date-cell formatting from the preceding (template) row so the new
# H/I cells reuse the existing date style instead of Excel minting new ones.
$ws.Range("H218:I218").Copy()
$ws.Range("H219:I219").PasteSpecial(-4122)
$ws.Range("H218:I218").Copy()
$ws.Range("H220:I220").PasteSpecial(-4122)
$ws.Range("H218:I218").Copy()
$ws.Range("H221:I221").PasteSpecial(-4122)

# Update active selection to match post-edit state
$ws.Range("F215").Select()

Write-Output "done"
